$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 2451200
$ws.Range("E8").Value = 2171400
$ws.Range("F8").Value = 1965700
$ws.Range("G8").Value = 1833700
$ws.Range("H8").Value = 1804100
$ws.Range("I8").Value = 1848200
$ws.Range("J8").Value = 1921400
$ws.Range("D9").Value = 1281400
$ws.Range("E9").Value = 1079700
$ws.Range("F9").Value = 919300
$ws.Range("G9").Value = 828700
$ws.Range("H9").Value = 790400
$ws.Range("I9").Value = 853600
$ws.Range("J9").Value = 942300
$ws.Range("D10").Value = 1169800
$ws.Range("E10").Value = 1091700
$ws.Range("F10").Value = 1046500
$ws.Range("G10").Value = 1005000
$ws.Range("H10").Value = 1013800
$ws.Range("I10").Value = 994600
$ws.Range("J10").Value = 979100
$ws.Range("G12").Value = 4800
$ws.Range("E14").Value = -29700
$ws.Range("D15").Value = 158300
$ws.Range("E15").Value = 154300
$ws.Range("F15").Value = 148200
$ws.Range("G15").Value = 143200
$ws.Range("H15").Value = 140000
$ws.Range("I15").Value = 134100
$ws.Range("J15").Value = 130300
$ws.Range("D17").Value = 1875600
$ws.Range("E17").Value = 1662300
$ws.Range("F17").Value = 1422500
$ws.Range("G17").Value = 1436900
$ws.Range("H17").Value = 1291400
$ws.Range("I17").Value = 1505200
$ws.Range("J17").Value = 1361200
$ws.Range("D18").Value = 575600
$ws.Range("E18").Value = 509200
$ws.Range("F18").Value = 543300
$ws.Range("G18").Value = 396800
$ws.Range("H18").Value = 512800
$ws.Range("I18").Value = 343000
$ws.Range("J18").Value = 560300
$ws.Range("D20").Value = 313200
$ws.Range("E20").Value = 247300
$ws.Range("F20").Value = 274000
$ws.Range("G20").Value = 224200
$ws.Range("H20").Value = 274300
$ws.Range("I20").Value = 276900
$ws.Range("J20").Value = 225500
$ws.Range("D21").Value = 1047100
$ws.Range("E21").Value = 762500
$ws.Range("F21").Value = 965500
$ws.Range("G21").Value = 624200
$ws.Range("H21").Value = 927000
$ws.Range("I21").Value = 623700
$ws.Range("J21").Value = 916100
$ws.Range("D22").Value = 77600
$ws.Range("E22").Value = 81800
$ws.Range("F22").Value = 78400
$ws.Range("G22").Value = 80300
$ws.Range("H22").Value = 73500
$ws.Range("I22").Value = 81300
$ws.Range("J22").Value = 62500
$ws.Range("D23").Value = 811200
$ws.Range("E23").Value = 674700
$ws.Range("F23").Value = 738900
$ws.Range("G23").Value = 540700
$ws.Range("H23").Value = 713600
$ws.Range("I23").Value = 538700
$ws.Range("J23").Value = 723300
$ws.Range("D24").Value = 126900
$ws.Range("E24").Value = 119300
$ws.Range("F24").Value = 103600
$ws.Range("G24").Value = 102600
$ws.Range("H24").Value = 98200
$ws.Range("I24").Value = 98600
$ws.Range("J24").Value = 121400
$ws.Range("D26").Value = 684200
$ws.Range("E26").Value = 555400
$ws.Range("F26").Value = 635300
$ws.Range("G26").Value = 438100
$ws.Range("H26").Value = 615400
$ws.Range("I26").Value = 440100
$ws.Range("J26").Value = 601900
$ws.Range("D27").Value = 610100
$ws.Range("E27").Value = 478100
$ws.Range("F27").Value = 569700
$ws.Range("G27").Value = 383400
$ws.Range("H27").Value = 551700
$ws.Range("I27").Value = 395500
$ws.Range("J27").Value = 534700
$ws.Range("D32").Value = -313200
$ws.Range("E32").Value = -247300
$ws.Range("F32").Value = -274000
$ws.Range("G32").Value = -224200
$ws.Range("H32").Value = -274300
$ws.Range("I32").Value = -276900
$ws.Range("J32").Value = -225500
$ws.Range("D33").Value = 610100
$ws.Range("E33").Value = 478100
$ws.Range("F33").Value = 569700
$ws.Range("G33").Value = 383400
$ws.Range("H33").Value = 551700
$ws.Range("I33").Value = 395500
$ws.Range("J33").Value = 534700
$ws.Range("D35").Value = 610100
$ws.Range("E35").Value = 478100
$ws.Range("F35").Value = 569700
$ws.Range("G35").Value = 383400
$ws.Range("H35").Value = 551700
$ws.Range("I35").Value = 395500
$ws.Range("J35").Value = 534700
$ws.Range("D41").Value = 1502900
$ws.Range("E41").Value = 1370500
$ws.Range("F41").Value = 1113100
$ws.Range("G41").Value = 1028800
$ws.Range("H41").Value = 906800
$ws.Range("I41").Value = 1519200
$ws.Range("J41").Value = 1751800
$ws.Range("D42").Value = 217700
$ws.Range("E42").Value = 284400
$ws.Range("F42").Value = 244000
$ws.Range("G42").Value = 450400
$ws.Range("H42").Value = 700400
$ws.Range("I42").Value = 170600
$ws.Range("J42").Value = 56000
$ws.Range("D43").Value = 972300
$ws.Range("E43").Value = 909200
$ws.Range("F43").Value = 810800
$ws.Range("G43").Value = 788400
$ws.Range("H43").Value = 813300
$ws.Range("I43").Value = 870600
$ws.Range("J43").Value = 841500
$ws.Range("D44").Value = 305900
$ws.Range("E44").Value = 328400
$ws.Range("F44").Value = 286200
$ws.Range("G44").Value = 268800
$ws.Range("H44").Value = 268800
$ws.Range("I44").Value = 291900
$ws.Range("J44").Value = 296100
$ws.Range("D45").Value = 192000
$ws.Range("E45").Value = 211300
$ws.Range("F45").Value = 178400
$ws.Range("G45").Value = 160400
$ws.Range("H45").Value = 163800
$ws.Range("I45").Value = 158300
$ws.Range("J45").Value = 237000
$ws.Range("D46").Value = 3190900
$ws.Range("E46").Value = 3104000
$ws.Range("F46").Value = 2632500
$ws.Range("G46").Value = 2697000
$ws.Range("H46").Value = 2853000
$ws.Range("I46").Value = 3010600
$ws.Range("J46").Value = 3182300
$ws.Range("D47").Value = 5150200
$ws.Range("E47").Value = 4948100
$ws.Range("F47").Value = 4790200
$ws.Range("G47").Value = 4482200
$ws.Range("H47").Value = 4436000
$ws.Range("I47").Value = 4281400
$ws.Range("J47").Value = 4318500
$ws.Range("D48").Value = 7402800
$ws.Range("E48").Value = 7209200
$ws.Range("F48").Value = 6703700
$ws.Range("G48").Value = 6361700
$ws.Range("H48").Value = 6269600
$ws.Range("I48").Value = 6136200
$ws.Range("J48").Value = 6059800
$ws.Range("D49").Value = 742800
$ws.Range("E49").Value = 749500
$ws.Range("F49").Value = 721700
$ws.Range("G49").Value = 709900
$ws.Range("H49").Value = 727200
$ws.Range("I49").Value = 741300
$ws.Range("J49").Value = 762600
$ws.Range("D52").Value = 699400
$ws.Range("E52").Value = 685200
$ws.Range("F52").Value = 645100
$ws.Range("G52").Value = 621300
$ws.Range("H52").Value = 621600
$ws.Range("I52").Value = 552100
$ws.Range("J52").Value = 535900
$ws.Range("D54").Value = 17186100
$ws.Range("E54").Value = 16696000
$ws.Range("F54").Value = 15493300
$ws.Range("G54").Value = 14871900
$ws.Range("H54").Value = 14907300
$ws.Range("I54").Value = 14721600
$ws.Range("J54").Value = 14859100
$ws.Range("D57").Value = 1023000
$ws.Range("E57").Value = 1962800
$ws.Range("F57").Value = 1664300
$ws.Range("G57").Value = 1637400
$ws.Range("H57").Value = 1547500
$ws.Range("I57").Value = 1593500
$ws.Range("J57").Value = 1528300
$ws.Range("D58").Value = 2209700
$ws.Range("E58").Value = 2007300
$ws.Range("F58").Value = 973600
$ws.Range("G58").Value = 758200
$ws.Range("H58").Value = 824300
$ws.Range("I58").Value = 1237200
$ws.Range("J58").Value = 1236000
$ws.Range("D59").Value = 956300
$ws.Range("E59").Value = 99800
$ws.Range("F59").Value = 115700
$ws.Range("G59").Value = 94600
$ws.Range("H59").Value = 125300
$ws.Range("I59").Value = 122200
$ws.Range("J59").Value = 159600
$ws.Range("D60").Value = 4188900
$ws.Range("E60").Value = 4069900
$ws.Range("F60").Value = 2753500
$ws.Range("G60").Value = 2490200
$ws.Range("H60").Value = 2497100
$ws.Range("I60").Value = 2953000
$ws.Range("J60").Value = 2923900
$ws.Range("D61").Value = 2880900
$ws.Range("E61").Value = 2695800
$ws.Range("F61").Value = 3426500
$ws.Range("G61").Value = 3477300
$ws.Range("H61").Value = 3518000
$ws.Range("I61").Value = 2976300
$ws.Range("J61").Value = 2963300
$ws.Range("D62").Value = 994700
$ws.Range("E62").Value = 981700
$ws.Range("F62").Value = 934700
$ws.Range("G62").Value = 892000
$ws.Range("H62").Value = 871000
$ws.Range("I62").Value = 879800
$ws.Range("J62").Value = 867800
$ws.Range("D66").Value = 9047800
$ws.Range("E66").Value = 8696800
$ws.Range("F66").Value = 7987700
$ws.Range("G66").Value = 7701700
$ws.Range("H66").Value = 7759500
$ws.Range("I66").Value = 7690300
$ws.Range("J66").Value = 7668400
$ws.Range("D72").Value = 6987400
$ws.Range("E72").Value = 6772900
$ws.Range("F72").Value = 6491200
$ws.Range("G72").Value = 6299200
$ws.Range("H72").Value = 6118700
$ws.Range("I72").Value = 5906400
$ws.Range("J72").Value = 5746800
$ws.Range("D76").Value = 8138300
$ws.Range("E76").Value = 7999200
$ws.Range("F76").Value = 7505600
$ws.Range("G76").Value = 7170300
$ws.Range("H76").Value = 7147900
$ws.Range("I76").Value = 7031300
$ws.Range("J76").Value = 7190700
$ws.Range("D81").Value = 610100
$ws.Range("E81").Value = 478100
$ws.Range("F81").Value = 569700
$ws.Range("G81").Value = 383400
$ws.Range("H81").Value = 551700
$ws.Range("I81").Value = 395500
$ws.Range("J81").Value = 534700
$ws.Range("D89").Value = 638300
$ws.Range("E89").Value = 526300
$ws.Range("F89").Value = 559600
$ws.Range("G89").Value = 521300
$ws.Range("H89").Value = 557600
$ws.Range("I89").Value = 529900
$ws.Range("J89").Value = 524500
$ws.Range("D91").Value = -411900
$ws.Range("E91").Value = -450500
$ws.Range("F91").Value = -314400
$ws.Range("G91").Value = -435600
$ws.Range("H91").Value = -336000
$ws.Range("I91").Value = -419600
$ws.Range("J91").Value = -348400
$ws.Range("D94").Value = -354800
$ws.Range("E94").Value = -148900
$ws.Range("F94").Value = -63900
$ws.Range("H94").Value = -812500
$ws.Range("I94").Value = -423300
$ws.Range("J94").Value = -180400
$ws.Range("D96").Value = -409800
$ws.Range("E96").Value = -213800
$ws.Range("F96").Value = -372600
$ws.Range("G96").Value = -194400
$ws.Range("H96").Value = -338700
$ws.Range("I96").Value = -176800
$ws.Range("J96").Value = -308000
$ws.Range("D100").Value = -148100
$ws.Range("E100").Value = -132800
$ws.Range("F100").Value = -420000
$ws.Range("G100").Value = -371400
$ws.Range("H100").Value = -351300
$ws.Range("I100").Value = -305700
$ws.Range("J100").Value = -200900
$ws.Range("E101").Value = 12800
$ws.Range("F101").Value = 8600
$ws.Range("G101").Value = -16400
$ws.Range("I101").Value = -33400
$ws.Range("D102").Value = 132400
$ws.Range("E102").Value = 257400
$ws.Range("F102").Value = 84300
$ws.Range("G102").Value = 122000
$ws.Range("H102").Value = -612500
$ws.Range("I102").Value = -232600
$ws.Range("J102").Value = 146000
